$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (rows 2:11 shift down to 3:12), pushing down all
# the existing merged ranges with it.
$ws.Rows("2:2").Insert()

# Make the new row look like the rest of the data rows (centered, regular
# weight) rather than the bold header-style formatting that Insert() copied
# down from row 1.
$ws.Range("A2:L2").Font.Bold = $false

# Populate the new row with the new model entry.
$ws.Range("A2").Value = "Shorthouse_etal_ionchannels_base.json"
$ws.Range("C2").Value = "Model of core osmotic network in resting state"

# Re-create the merged layout for the new row, matching every other data row.
$ws.Range("A2:B2").Merge()
$ws.Range("C2:L2").Merge()

# Match the saved selection state of the edited workbook.
$ws.Range("C9:L9").Select() | Out-Null

# Restore the explicit portrait orientation recorded on the sheet.
$ws.PageSetup.Orientation = 1
